$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Carlos Cadiz"
$ws.Range("B11").Value = "carloscadiz2805@gmail.com"
$ws.Range("C11").Value = "21019653-6"
$ws.Range("D11").Value = "ALUMNO"
$ws.Range("E11").Value = "Ingeniería de Ejecución en Computación"
